$p = $ppt.ActivePresentation

# --- Slide 1: "Correlatore: ... Marco   LANZA" ---
# Collapse the trailing " " + "LANZA" runs into a single run " LANZA"
# (keeping the rPr of the run that carried the first of those two runs).
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$run1 = $tr1.Characters(51, 6)
$run1.Text = " LANZA"

# --- Slide 10: "Foto scattata durante l'Arduino Day (16/05/2019)" ---
# Collapse "l’" + "Arduino" into a single run "l’Arduino"
# Collapse " " + "Day (16/05/2019)" into a single run " Day (16/05/2019)"
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(3)
$tr10 = $sh10.TextFrame.TextRange
$run2 = $tr10.Characters(23, 9)
$run2.Text = "l’Arduino"
$run3 = $tr10.Characters(32, 17)
$run3.Text = " Day (16/05/2019)"
